$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "iaest-dimension:ccaa-nombre"
$ws.Range("B2").Value = "iaest-measure:medidas-aumento-valor-anadido"
$ws.Range("C2").Value = "iaest-measure:ayudas-natura-2000"
$ws.Range("D2").Value = "iaest-measure:pagos-otras-ayudas"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "iaest-measure:medidas-modernizacion"
$ws.Range("G2").Value = "iaest-measure:diversificacion-actividades-no-agricolas"
$ws.Range("H2").Value = "null"
$ws.Range("I2").Value = "null"
$ws.Range("J2").Value = "null"
$ws.Range("K2").Value = "iaest-measure:medidas-asesoramiento"
$ws.Range("L2").Value = "sdmx-dimension:refArea"
$ws.Range("M2").Value = "iaest-measure:pagos-agricultura-ecologica"
$ws.Range("N2").Value = "iaest-measure:ayudas-relativas-bienestar-animales"
$ws.Range("O2").Value = "iaest-measure:medidas-marco-del-agua"
$ws.Range("P2").Value = "iaest-measure:medidas-cumplimiento-normas"
$ws.Range("Q2").Value = "iaest-measure:fomento-actividades-turisticas"
$ws.Range("R2").Value = "iaest-measure:medidas-participacion-programas-calidad"
$ws.Range("S2").Value = "iaest-measure:identificador"
$ws.Range("T2").Value = "sdmx-dimension:refArea"
$ws.Range("A3").Value = "dim"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "null"
$ws.Range("I3").Value = "null"
$ws.Range("J3").Value = "null"
$ws.Range("K3").Value = "medida"
$ws.Range("L3").Value = "dim"
$ws.Range("M3").Value = "medida"
$ws.Range("N3").Value = "medida"
$ws.Range("O3").Value = "medida"
$ws.Range("P3").Value = "medida"
$ws.Range("Q3").Value = "medida"
$ws.Range("R3").Value = "medida"
$ws.Range("S3").Value = "medida"
$ws.Range("T3").Value = "dim"
$ws.Range("A4").Value = "skos:Concept"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "null"
$ws.Range("I4").Value = "null"
$ws.Range("J4").Value = "null"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("L4").Value = "URI-Municipio"
$ws.Range("M4").Value = "xsd:int"
$ws.Range("N4").Value = "xsd:int"
$ws.Range("O4").Value = "xsd:int"
$ws.Range("P4").Value = "xsd:int"
$ws.Range("Q4").Value = "xsd:int"
$ws.Range("R4").Value = "xsd:int"
$ws.Range("S4").Value = "xsd:int"
$ws.Range("T4").Value = "URI-comarca"
$ws.Range("B5").Clear()
$ws.Range("C5").Clear()
$ws.Range("D5").Clear()
$ws.Range("F5").Clear()
$ws.Range("G5").Clear()
$ws.Range("K5").Clear()
$ws.Range("M5").Clear()
$ws.Range("N5").Clear()
$ws.Range("O5").Clear()
$ws.Range("P5").Clear()
$ws.Range("Q5").Clear()
$ws.Range("R5").Clear()
